$wb = $excel.ActiveWorkbook
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws4 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws4.Name = "Sheet4"
Write-Host $wb.Worksheets.Count
for ($i = 1; $i -le $wb.Worksheets.Count; $i++) {
    Write-Host $wb.Worksheets.Item($i).Name
}
